$d = $word.ActiveDocument

# Bounds is a list of offsets (relative to BaseStart) that mark the
# boundaries between the desired runs, e.g. @(0, 5, 10) creates a
# split between [0,5) and [5,10). Toggling a character formatting
# property on/off over each sub-range forces the engine to keep that
# sub-range as its own run, without altering the visible formatting.
function Split-RunsAtOffsets($BaseStart, $Bounds) {
    for ($i = 1; $i -lt $Bounds.Length; $i++) {
        $s = $BaseStart + $Bounds[$i - 1]
        $e = $BaseStart + $Bounds[$i]
        if ($e -gt $s) {
            $p = $d.Range($s, $e)
            $p.Font.Bold = 1
            $p.Font.Bold = 0
        }
    }
}

# ---------------------------------------------------------------------
# 1) "***** 2020" -> "July ***** 2020"  (new run "July " before it)
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("***** 2020", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$ins = $d.Range($start, $start)
$ins.InsertBefore("July ")
Split-RunsAtOffsets $start @(0, 5)

# ---------------------------------------------------------------------
# 2) ". This study focussed on 2016, which had the largest spatial
#    coverage of sampling stations in an effort to resolve the primary
#    migration pathways through the region. The previous year 2015 also
#    had similar spatial coverage of sampling but there was expected and
#    observed lower pink abundance, due to their biennial life patterns."
#    ->
#    ". This study focussed on 2015 and 2016, which had the highest
#    frequency of sampling stations in an effort to resolve the primary
#    migration dynamics."
# ---------------------------------------------------------------------
$oldText2 = "2016, which had the largest spatial coverage of sampling stations in an effort to resolve the primary migration pathways through the region. The previous year 2015 also had similar spatial coverage of sampling but there was expected and observed lower pink abundance, due to their biennial life patterns."
$newText2 = "2015 and 2016, which had the highest frequency of sampling stations in an effort to resolve the primary migration dynamics."
$rng2 = $d.Content
$rng2.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, $newText2, 2)

$fullText2 = ". This study focussed on 2015 and 2016, which had the highest frequency of sampling stations in an effort to resolve the primary migration dynamics."
$find2 = $d.Content
$find2.Find.Execute($fullText2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base2 = $find2.Start
# part boundaries (character offsets from $base2):
#   0   ". This study focussed on "
#   25  "2015 and "
#   34  "2016, which had the "
#   54  "highest frequency of"
#   74  " sampling stations in an effort to resolve the primary migration "
#   139 "dynamics"
#   147 "."
#   148
Split-RunsAtOffsets $base2 @(0, 25, 34, 54, 74, 139, 147, 148)

# ---------------------------------------------------------------------
# 3) "For this study, six sites (three from each region) were selected"
#    -> "For this study, two sites (one from each region) were selected"
# ---------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("six sites (three from each region)", $true, $false, $false, $false, $false, $true, 1, $false, "two sites (one from each region)", 2)

$fullText3 = ". Sites were sampled every 4-7 days throughout the season, depending on weather conditions. For this study, two sites (one from each region) were selected, in order to obtain a sample size of 10 pink and 10 chum per set (n=120 total), still acquiring sufficient coverage for each region. The dates were chosen in mid-June (Table 1) to align with the peak out-migration of salmon "
$find3 = $d.Content
$find3.Find.Execute($fullText3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base3 = $find3.Start
# part boundaries (character offsets from $base3):
#   0   ". Sites were sampled every 4-7 days throughout the season, depending on weather conditions. For this study, "
#   108 "two"
#   111 " sites ("
#   119 "one"
#   122 " from each region) were selected, in order to obtain a sample size of 10 pink and 10 chum per set (n=120 total), still acquiring sufficient coverage for each region. The dates were chosen in mid-June (Table 1) to align with the peak out-migration of salmon "
#   379
Split-RunsAtOffsets $base3 @(0, 108, 111, 119, 122, 379)

# ---------------------------------------------------------------------
# 4) Merge the "Microplastics were not the focus..." runs into a single
#    run (formatting-only change; visible text is unchanged). Re-typing
#    the text via Find/Replace over the full span naturally merges the
#    previously split runs (which all share identical formatting) into
#    one, while leaving the preceding <w:tab/> run untouched. The engine
#    also tends to sweep any immediately-following same-format runs into
#    that merge, so force a split right after the replaced span, and
#    restore the original run boundaries of the (untouched) text that
#    follows it.
# ---------------------------------------------------------------------
$oldText4 = "Microplastics were not the focus of this study but they were found in *****% of juvenile salmon stomachs, and one macroplastic was found to be 30% weight of a pink salmon stomach. That 6 mm macroplastic had the shape, color and texture of a broken straw piece and appeared larger than the sphincter could potentially pass, which would likely reduce survival for that fish. Impacts of plastics on salmon and occurrence in empty stomachs"
$rng4 = $d.Content
$rng4.Find.Execute($oldText4, $true, $false, $false, $false, $false, $true, 1, $false, $oldText4, 2)

$find4 = $d.Content
$find4.Find.Execute($oldText4, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s4 = $find4.Start
$e4 = $find4.End
Split-RunsAtOffsets $s4 @(0, ($e4 - $s4))
# restore original run boundaries for the untouched trailing text:
# "," / " with potential for cumulative effects, " / "should be researched..."
Split-RunsAtOffsets $e4 @(0, 1, 41, 124)
